$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "65.919.74", "  +0.72%  ") are preserved as text, matching the source data.
$ws.Range('D2:E51').NumberFormat = '@'

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range('D2').Value = '65.919.74'
$ws.Range('E2').Value = '  +0.72%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range('D3').Value = '3.305.60'
$ws.Range('E3').Value = '  -0.04%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range('E4').Value = '  +0.16%  '

# Row 5: 'Solana' -> 'Solana'
$ws.Range('D5').Value = '187.54'
$ws.Range('E5').Value = '  +4.62%  '

# Row 6: 'BNB' -> 'BNB'
$ws.Range('D6').Value = '555.73'
$ws.Range('E6').Value = '  +0.04%  '

# Row 7: 'USDC' -> 'USDC'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.03%  '

# Row 8: 'XRP' -> 'XRP'
$ws.Range('D8').Value = '0.582'
$ws.Range('E8').Value = '  -0.78%  '

# Row 9: 'LidoStakedEther' -> 'LidoStakedEther'
$ws.Range('D9').Value = '3.301.45'
$ws.Range('E9').Value = '  -0.04%  '

# Row 10: 'Dogecoin' -> 'Dogecoin'
$ws.Range('D10').Value = '0.183'
$ws.Range('E10').Value = '  -0.34%  '

# Row 11: 'Cardano' -> 'Cardano'
$ws.Range('D11').Value = '0.584'
$ws.Range('E11').Value = '  +0.55%  '

# Row 12: 'Avalanche' -> 'Avalanche'
$ws.Range('D12').Value = '47.24'
$ws.Range('E12').Value = '  +0.46%  '

# Row 13: 'ShibaInu' -> 'ShibaInu'
$ws.Range('D13').Value = '0.0000269'
$ws.Range('E13').Value = '  +2.57%  '

# Row 14: 'Polkadot' -> 'Polkadot'
$ws.Range('D14').Value = '8.67'
$ws.Range('E14').Value = '  +2.16%  '

# Row 15: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range('D15').Value = '3.839.75'
$ws.Range('E15').Value = '  +0.21%  '

# Row 16: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range('D16').Value = '603.91'
$ws.Range('E16').Value = '  +1.09%  '

# Row 17: 'Chainlink' -> 'WrappedBTC'
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '65.964.57'
$ws.Range('E17').Value = '  +0.83%  '

# Row 18: 'WrappedBTC' -> 'Chainlink'
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '17.96'
$ws.Range('E18').Value = '  -0.17%  '

# Row 19: 'TRON' -> 'TRON'
$ws.Range('E19').Value = '  +1.12%  '

# Row 20: 'WrappedEther' -> 'WrappedEther'
$ws.Range('D20').Value = '3.283.01'
$ws.Range('E20').Value = '  -0.25%  '

# Row 21: 'Uniswap' -> 'Uniswap'
$ws.Range('D21').Value = '11.06'
$ws.Range('E21').Value = '  -2.54%  '

# Row 22: 'Polygon' -> 'Polygon'
$ws.Range('D22').Value = '0.906'
$ws.Range('E22').Value = '  +1.11%  '

# Row 23: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range('D23').Value = '18.33'
$ws.Range('E23').Value = '  +7.70%  '

# Row 24: 'Toncoin' -> 'Toncoin'
$ws.Range('D24').Value = '5.05'
$ws.Range('E24').Value = '  +0.30%  '

# Row 25: 'Litecoin' -> 'Litecoin'
$ws.Range('D25').Value = '100.06'
$ws.Range('E25').Value = '  -1.70%  '

# Row 26: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range('D26').Value = '3.96'
$ws.Range('E26').Value = '  +0.00%  '

# Row 27: 'ImmutableX' -> 'ImmutableX'
$ws.Range('D27').Value = '2.76'
$ws.Range('E27').Value = '  +4.10%  '

# Row 28: 'LEO' -> 'LEO'
$ws.Range('D28').Value = '5.93'
$ws.Range('E28').Value = '  -0.81%  '

# Row 29: 'RenderToken' -> 'RenderToken'
$ws.Range('D29').Value = '9.56'
$ws.Range('E29').Value = '  +3.57%  '

# Row 30: 'Filecoin' -> 'Filecoin'
$ws.Range('D30').Value = '8.69'
$ws.Range('E30').Value = '  +0.95%  '

# Row 31: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range('D31').Value = '30.27'
$ws.Range('E31').Value = '  -0.53%  '

# Row 32: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range('D32').Value = '6.72'
$ws.Range('E32').Value = '  +8.56%  '

# Row 33: 'dogwifhat' -> 'dogwifhat'
$ws.Range('E33').Value = '  +1.62%  '

# Row 34: 'Bittensor' -> 'Bittensor'
$ws.Range('D34').Value = '579.69'
$ws.Range('E34').Value = '  +11.83%  '

# Row 35: 'Cosmos' -> 'Cosmos'
$ws.Range('D35').Value = '11.07'
$ws.Range('E35').Value = '  +0.98%  '

# Row 36: 'Hedera' -> 'Hedera'
$ws.Range('E36').Value = '  +0.59%  '

# Row 37: 'Dai' -> 'Maker'
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '3.704.35'
$ws.Range('E37').Value = '  -2.22%  '

# Row 38: 'OKB' -> 'Dai'
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.13%  '

# Row 39: 'Maker' -> 'OKB'
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '57.03'
$ws.Range('E39').Value = '  +0.76%  '

# Row 40: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range('D40').Value = '33.91'
$ws.Range('E40').Value = '  +6.49%  '

# Row 41: 'PEPE' -> 'PEPE'
$ws.Range('D41').Value = '0.0₃0717'
$ws.Range('E41').Value = '  +1.47%  '

# Row 42: 'Kaspa' -> 'Kaspa'
$ws.Range('E42').Value = '  +5.50%  '

# Row 43: 'CoreDAO' -> 'Stacks'
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '3.28'
$ws.Range('E43').Value = '  -4.60%  '

# Row 44: 'Stacks' -> 'CoreDAO'
$ws.Range('B44').Value = 'CoreDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D44').Value = '3.41'
$ws.Range('E44').Value = '  +4.00%  '

# Row 45: 'Fetch.AI' -> 'Fetch.AI'
$ws.Range('D45').Value = '2.67'
$ws.Range('E45').Value = '  +1.42%  '

# Row 46: 'TheGraph' -> 'TheGraph'
$ws.Range('D46').Value = '0.340'
$ws.Range('E46').Value = '  +0.98%  '

# Row 47: 'ApeXProtocol' -> 'ApeXProtocol'
$ws.Range('D47').Value = '3.38'
$ws.Range('E47').Value = '  +2.79%  '

# Row 48: 'VeChain' -> 'VeChain'
$ws.Range('D48').Value = '0.0420'
$ws.Range('E48').Value = '  +2.69%  '

# Row 49: 'Stellar' -> 'Stellar'
$ws.Range('E49').Value = '  +0.39%  '

# Row 50: 'ThetaToken' -> 'ThetaToken'
$ws.Range('E50').Value = '  -0.08%  '

# Row 51: 'FirstDigitalUSD' -> 'FirstDigitalUSD'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.20%  '

# Restore default (General) formatting so the cells carry no explicit style,
# matching the original workbook layout.
$ws.Range('D2:E51').ClearFormats()